$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 45048.84104810185

$ws.Range("A4").Value = "2/05/2023"
$ws.Range("B4").Value = "Волоьбуев"
$ws.Range("C4").Value = "8935910611"
$ws.Range("D4").Value = "Cordiant"
$ws.Range("E4").Value = "Gravity"
$ws.Range("F4").Value = "195/60 R17"
$ws.Range("G4").Value = "6684"
$ws.Range("H4").Value = "6584168"
$ws.Range("I4").Value = "681"
$ws.Range("J4").Value = "Диалло"
$ws.Range("K4").Value = 45048.94851479166
$ws.Range("K4").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("L4").Value = "COR812619"

$ws.Range("A5").Value = "2/0506/65146"
$ws.Range("B5").Value = "вололуев"
$ws.Range("C5").Value = "89535919844"
$ws.Range("D5").Value = "Tunga"
$ws.Range("E5").Value = "Zodiak 2"
$ws.Range("F5").Value = "185/45 R19"
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = "234235"
$ws.Range("J5").Value = "Никифоров"
$ws.Range("K5").Value = 45048.95025267241
$ws.Range("K5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("L5").Value = "TUN809835"
